$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-16 (columns A-T), reflecting the updated
# TPM-derived NATMI ligand-receptor computation. "ECs" is now also a sending
# cluster (in addition to being a target cluster), yielding 15 data rows
# (5 sending clusters x 3 target clusters) instead of the original 12.
$data = @(
    @("ECs", "Tnfsf8", "Tnfrsf8", "ECs", 2, 0.6666666666666666, 0.06805433333333334, 0.204163, 0.01286726242426351, 0.01295654885662802, 3, 1, 0.4756293333333333, 1.426888, 0.4551485584980385, 0.4551485584980385, 0.03236863719377778, 0.291317734744, 0.005856515944219514, 0.005897154535203652),
    @("ECs", "Tnfsf8", "Tnfrsf8", "Inflammatory-Mac", 2, 0.6666666666666666, 0.06805433333333334, 0.204163, 0.01286726242426351, 0.01295654885662802, 1, 0.3333333333333333, 0.3753246666666667, 1.125974, 0.3591630478399639, 0.3591630478399639, 0.02554246997355556, 0.229882229762, 0.004621445189655126, 0.004653513576833919),
    @("ECs", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 2, 0.6666666666666666, 0.06805433333333334, 0.204163, 0.01286726242426351, 0.01295654885662802, 1, 0.3333333333333333, 0.194044, 0.582132, 0.1856883936619974, 0.1856883936619974, 0.01320553505733334, 0.118849815516, 0.00238930129038887, 0.002405880744590446),
    @("FAPs", "Tnfsf8", "Tnfrsf8", "ECs", 3, 1, 1.400976666666667, 4.20293, 0.2648873853774182, 0.2667254492047413, 3, 1, 0.4756293333333333, 1.426888, 0.4551485584980385, 0.4551485584980385, 0.6663455979822223, 5.99711038184, 0.1205631116188463, 0.1213997037202798),
    @("FAPs", "Tnfsf8", "Tnfrsf8", "Inflammatory-Mac", 3, 1, 1.400976666666667, 4.20293, 0.2648873853774182, 0.2667254492047413, 1, 0.3333333333333333, 0.3753246666666667, 1.125974, 0.3591630478399639, 0.3591630478399639, 0.5258211004244445, 4.732389903820001, 0.09513776066651261, 0.09579792527285838),
    @("FAPs", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 3, 1, 1.400976666666667, 4.20293, 0.2648873853774182, 0.2667254492047413, 1, 0.3333333333333333, 0.194044, 0.582132, 0.1856883936619974, 0.1856883936619974, 0.2718511163066667, 2.44666004676, 0.04918651309205924, 0.0495278202116031),
    @("Inflammatory-Mac", "Tnfsf8", "Tnfrsf8", "ECs", 3, 1, 2.224606, 6.673818000000001, 0.4206137624240114, 0.4235324176135906, 3, 1, 0.4756293333333333, 1.426888, 0.4551485584980385, 0.4551485584980385, 1.058087868709333, 9.522790818384001, 0.1914417476517253, 0.192770169354015),
    @("Inflammatory-Mac", "Tnfsf8", "Tnfrsf8", "Inflammatory-Mac", 3, 1, 2.224606, 6.673818000000001, 0.4206137624240114, 0.4235324176135906, 1, 0.3333333333333333, 0.3753246666666667, 1.125974, 0.3591630478399639, 0.3591630478399639, 0.8349495054146667, 7.514545548732001, 0.1510689208756424, 0.1521171939691256),
    @("Inflammatory-Mac", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 3, 1, 2.224606, 6.673818000000001, 0.4206137624240114, 0.4235324176135906, 1, 0.3333333333333333, 0.194044, 0.582132, 0.1856883936619974, 0.1856883936619974, 0.431671446664, 3.885043019976, 0.07810309389664369, 0.07864505429044991),
    @("MuSCs", "Tnfsf8", "Tnfrsf8", "ECs", 1, 0.5, 0.109342, 0.218684, 0.02067366086892072, 0.01387807746831131, 3, 1, 0.4756293333333333, 1.426888, 0.4551485584980385, 0.4551485584980385, 0.05200626256533333, 0.3120375753919999, 0.009409586943366573, 0.006316586954426),
    @("MuSCs", "Tnfsf8", "Tnfrsf8", "Inflammatory-Mac", 1, 0.5, 0.109342, 0.218684, 0.02067366086892072, 0.01387807746831131, 1, 0.3333333333333333, 0.3753246666666667, 1.125974, 0.3591630478399639, 0.3591630478399639, 0.04103874970266667, 0.246232498216, 0.007425215047691365, 0.00498449260167782),
    @("MuSCs", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 1, 0.5, 0.109342, 0.218684, 0.02067366086892072, 0.01387807746831131, 1, 0.3333333333333333, 0.194044, 0.582132, 0.1856883936619974, 0.1856883936619974, 0.021217159048, 0.127302954288, 0.003838858877862783, 0.002576997912207486),
    @("Resolving-Mac", "Tnfsf8", "Tnfrsf8", "ECs", 3, 1, 1.485973, 4.457919, 0.2809579289053862, 0.2829075068567288, 3, 1, 0.4756293333333333, 1.426888, 0.4551485584980385, 0.4551485584980385, 0.7067723473413333, 6.360951126072, 0.1278775963398809, 0.1287649439341141),
    @("Resolving-Mac", "Tnfsf8", "Tnfrsf8", "Inflammatory-Mac", 3, 1, 1.485973, 4.457919, 0.2809579289053862, 0.2829075068567288, 1, 0.3333333333333333, 0.3753246666666667, 1.125974, 0.3591630478399639, 0.3591630478399639, 0.5577223209006668, 5.019500888106001, 0.1009097060604624, 0.1016099224194682),
    @("Resolving-Mac", "Tnfsf8", "Tnfrsf8", "Resolving-Mac", 3, 1, 1.485973, 4.457919, 0.2809579289053862, 0.2829075068567288, 1, 0.3333333333333333, 0.194044, 0.582132, 0.1856883936619974, 0.1856883936619974, 0.288344144812, 2.595097303308, 0.05217062650504283, 0.0525326405031465)
)

$rowIdx = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $val
        $colIdx = $colIdx + 1
    }
    $rowIdx = $rowIdx + 1
}
